$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. date number format on column A) from the last existing row (48)
# down into the new rows before writing values, so the new date cells keep the
# "yyyy-mm-dd" style used throughout column A.
$ws.Range("A48:F48").Copy() | Out-Null
$ws.Range("A49:F60").PasteSpecial(-4122) | Out-Null

$data = @(
    @(44695, 0, 327254, 6345, 16, 0),
    @(44696, 0, 327271, 6345, 17, 0),
    @(44697, 0, 327284, 6345, 13, 0),
    @(44698, 0, 327298, 6346, 14, 1),
    @(44699, 0, 327314, 6346, 16, 0),
    @(44700, 0, 327335, 6346, 21, 0),
    @(44701, 0, 327349, 6346, 14, 0),
    @(44702, 0, 327363, 6346, 14, 0),
    @(44703, 0, 327373, 6346, 10, 0),
    @(44704, 0, 327379, 6346, 6, 0),
    @(44705, 0, 327397, 6347, 18, 1),
    @(44706, 0, 327411, 6347, 14, 0)
)

$startRow = 49
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Update the view to match the new scroll/selection position recorded in the workbook
$ws.Range("E58").Select() | Out-Null

Write-Host "Applied casos/obitos update through 2022-05-25"
